$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SanityTC")

# New header names added to row 1 (columns BC:BW), in left-to-right order.
$headerCols = @("BC","BD","BE","BF","BG","BH","BI","BJ","BK","BL","BM","BN","BO","BP","BQ","BR","BS","BT","BU","BV","BW")
$headerVals = @("imagechoice","reportingvalue","images","date","likedislike","ratingscale","dropdown","ranking","multitextbox","textbox1","textbox2","multidropdown","dropdown1","dropdown2","multiradio","QuestionOptions","multicheckbox","ratingradio","ratingdropdown","ratingscalegrid","matrixgrid")

# Set the new case-flag text in J3 first (matches shared-string insertion order).
$ws.Range("J3").Value = " Start with a Blank Survey Button not present on page."

# Populate the new header cells in row 1.
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $headerVals[$i]
}

# Extend the formatting of the new columns in rows 2 and 3 to match column BB
# (blank data cells under the new headers).
$ws.Range("BB1").Copy()
$ws.Range("BC1:BW1").PasteSpecial(-4122)

$ws.Range("BB2").Copy()
$ws.Range("BC2:BW2").PasteSpecial(-4122)

$ws.Range("BB3").Copy()
$ws.Range("BC3:BW3").PasteSpecial(-4122)

# Update the saved view: clear the scrolled-right position and move the
# active selection back to A3.
[void]$ws.Activate()
[void]$ws.Range("A3").Select()
